# ---------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# The previous "总计" (Grand-total) sheet (3rd tab) is repurposed to hold the
# new 2022-Q1 per-fund holdings detail (same shape as the existing
# "2021-Q3"/"2021-Q4" tabs), and a brand-new "总计" summary sheet is appended
# right after it, carrying over the old 2-row summary plus a new row for
# 2022-Q1.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# xlPasteFormats - used to clone an existing cell style (by index) onto a
# range without creating new style entries.
$xlPasteFormats = -4122

# -----------------------------------------------------------------------
# 1) Turn the old "总计" sheet into the "2022-Q1" per-fund detail sheet.
# -----------------------------------------------------------------------
$wsFund = $wb.Worksheets.Item(3)
$wsFund.Name = "2022-Q1"

# Grab the header style (bold/border/centered) that already lives on B1 so
# the three brand-new header cells (E1:H1) match it exactly.
$wsFund.Range("B1").Copy()
$wsFund.Range("E1:H1").PasteSpecial($xlPasteFormats)

$wsFund.Cells.Item(1, 2).Value = "基金代码"
$wsFund.Cells.Item(1, 3).Value = "基金名称"
$wsFund.Cells.Item(1, 4).Value = "基金规模"
$wsFund.Cells.Item(1, 5).Value = "股票总仓位"
$wsFund.Cells.Item(1, 6).Value = "仓位占比"
$wsFund.Cells.Item(1, 7).Value = "持有市值(亿元)"
$wsFund.Cells.Item(1, 8).Value = "仓位排名"

# Row data: code, name, fund size, stock position, position ratio,
# held market value (100M CNY), position rank.
$fundRows = @(
    @("519702","交银趋势优先混合","108.29","71.40","2.77","2.9996",2),
    @("001645","国泰大健康股票A","34.81","90.83","6.24","2.1721",8),
    @("020001","国泰金鹰增长灵活配置混合","17.38","91.77","8.84","1.5364",3),
    @("011251","华安聚嘉精选混合A","31.29","89.00","3.01","0.9418",2),
    @("011128","华安精致生活混合A","33.30","85.22","2.21","0.7359",7),
    @("160215","国泰价值经典灵活配置混合（LOF）","6.36","92.86","9.27","0.5896",2),
    @("011252","华安聚嘉精选混合C","14.89","89.00","3.01","0.4482",2),
    @("008370","国泰研究精选两年持有期混合","4.15","92.87","10.03","0.4162",3),
    @("009804","国泰研究优势混合","4.14","90.91","9.74","0.4032",2),
    @("160212","国泰估值优势混合 (LOF)","8.98","62.69","4.14","0.3718",6),
    @("011568","鹏华产业升级混合A","21.41","68.03","1.73","0.3704",10),
    @("000601","华宝创新优选混合","12.99","87.56","2.71","0.3520",9),
    @("010738","大成优选升级一年持有期混合A","3.79","89.02","8.11","0.3074",3),
    @("000006","西部利得量化成长混合A","21.21","88.21","1.28","0.2715",8),
    @("010779","西部利得量化优选一年持有期混合A","11.33","87.57","2.22","0.2515",9),
    @("011321","国泰大健康股票C","3.47","90.83","6.24","0.2165",8),
    @("011129","华安精致生活混合C","7.57","85.22","2.21","0.1673",7),
    @("008185","诺安研究优选混合","2.59","94.06","5.30","0.1373",6),
    @("002197","国泰鑫策略价值灵活配置混合","6.92","21.00","1.00","0.0692",2),
    @("001850","国泰安益灵活配置混合A","6.63","21.23","0.80","0.0530",3),
    @("000367","国泰安康定期支付混合A","5.28","21.92","0.86","0.0454",3),
    @("001922","国泰多策略收益灵活配置混合","6.86","24.58","0.66","0.0453",10),
    @("011331","鹏华远见成长混合型证券投资基金A","2.29","63.89","1.72","0.0394",9),
    @("010780","西部利得量化优选一年持有期混合C","1.44","87.57","2.22","0.0320",9),
    @("010834","国泰同益18个月持有期混合型证券投资基金A","2.23","34.49","1.26","0.0281",6),
    @("001242","博时中证淘金大数据100指数A","2.51","93.96","0.98","0.0246",6),
    @("001243","博时中证淘金大数据100指数I","2.51","93.96","0.98","0.0246",6),
    @("011228","西部利得量化成长混合C","1.59","88.21","1.28","0.0204",8),
    @("002061","国泰安康定期支付混合C","2.35","21.92","0.86","0.0202",3),
    @("004252","国泰安益灵活配置混合C","2.10","21.23","0.80","0.0168",3),
    @("010835","国泰同益18个月持有期混合型证券投资基金C","0.87","34.49","1.26","0.0110",6),
    @("010739","大成优选升级一年持有期混合C","0.09","89.02","8.11","0.0073",3),
    @("011569","鹏华产业升级混合C","0.34","68.03","1.73","0.0059",10),
    @("011332","鹏华远见成长混合型证券投资基金C","0.21","63.89","1.72","0.0036",9)
)

$fundRowCount = $fundRows.Length
$fundLastRow = $fundRowCount + 1

# Column A holds the 0-based row index and keeps the same bold/bordered
# style as the existing A2/A3 cells - clone that style down the column
# first, then fill in the values.
$wsFund.Range("A2").Copy()
$wsFund.Range("A2:A" + $fundLastRow).PasteSpecial($xlPasteFormats)

# Columns B:G are free-form text (fund code/name/size/position data all
# come through as strings, e.g. "001645" keeps its leading zero). Force
# text storage with a "@" number format, write the values, then strip the
# number-format override back off so the cells end up unstyled, matching
# the target sheet.
$wsFund.Range("B2:G" + $fundLastRow).NumberFormat = "@"

for ($i = 0; $i -lt $fundRowCount; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    $wsFund.Cells.Item($r, 1).Value = $i

    $wsFund.Cells.Item($r, 2).Value = $row[0]
    $wsFund.Cells.Item($r, 3).Value = $row[1]
    $wsFund.Cells.Item($r, 4).Value = $row[2]
    $wsFund.Cells.Item($r, 5).Value = $row[3]
    $wsFund.Cells.Item($r, 6).Value = $row[4]
    $wsFund.Cells.Item($r, 7).Value = $row[5]

    # Column H (position rank) is a real number.
    $wsFund.Cells.Item($r, 8).Value = $row[6]
}

$wsFund.Range("B2:G" + $fundLastRow).ClearFormats()

# -----------------------------------------------------------------------
# 2) Insert the new "总计" summary sheet right after "2022-Q1", carrying
#    the old 2-row summary plus a new leading row for 2022-Q1.
# -----------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Add($null, $wsFund)
$wsTotal.Name = "总计"

$wsFund.Range("B1").Copy()
$wsTotal.Range("B1:D1").PasteSpecial($xlPasteFormats)

$wsTotal.Cells.Item(1, 2).Value = "日期"
$wsTotal.Cells.Item(1, 3).Value = "持有数量(只)"
$wsTotal.Cells.Item(1, 4).Value = "持有市值(亿元)"

# date label, holdings count, held market value (100M CNY)
$totalRows = @(
    @("2022-Q1", 34, 13.14),
    @("2021-Q4", 15, 7.49),
    @("2021-Q3", 10, 5.72)
)

$wsFund.Range("A2").Copy()
$wsTotal.Range("A2:A4").PasteSpecial($xlPasteFormats)

$wsTotal.Range("B2:B4").NumberFormat = "@"

for ($i = 0; $i -lt $totalRows.Length; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]

    $wsTotal.Cells.Item($r, 1).Value = $i
    $wsTotal.Cells.Item($r, 2).Value = $row[0]
    $wsTotal.Cells.Item($r, 3).Value = $row[1]
    $wsTotal.Cells.Item($r, 4).Value = $row[2]
}

$wsTotal.Range("B2:B4").ClearFormats()

Write-Output "2022-Q1 detail + 总计 summary sheets written"
